$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-3 and extend the table with new rows 4-11
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'es'
$ws.Range("C2").Value = 'One term refers to multiple concepts'
$ws.Range("D2").Value = 'One scientific term has a different meaning depending on the context it is used in. e.g. volts and voltage in Physics. Use of the term kinetic energy in both Physics and in Biology.'
$ws.Range("E2").Value = 'Terminology'
$ws.Range("F2").Value = 'Problems with use of language and scientific terms, inconsistent and overlapping terminology.'
$ws.Rows.Item(2).RowHeight = 75

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 'es'
$ws.Range("C3").Value = 'One concept has many scientific names'
$ws.Range("D3").Value = 'Different terms are used to refer to the same concept. e.g. voltage is also referred to as potential difference. Confusion between voltage and charge.'
$ws.Range("E3").Value = 'Terminology'
$ws.Range("F3").Value = 'Problems with use of language and scientific terms, inconsistent and overlapping terminology.'
$ws.Rows.Item(3).RowHeight = 60

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'es'
$ws.Range("C4").Value = 'Scientific use of everyday language'
$ws.Range("D4").Value = 'Everyday terms that students reused in a scientific context, where their scientific meaning may be slightly different to that understood by students.e.g. in Physics, the “drop” part of “forward voltage drop”, "current" related to electricity and  "requency", relating to waves. Use of the word "proof" to mean evidence.'
$ws.Range("E4").Value = 'Terminology'
$ws.Range("F4").Value = 'Problems with use of language and scientific terms, inconsistent and overlapping terminology.'
$ws.Rows.Item(4).RowHeight = 120

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'es'
$ws.Range("C5").Value = 'Obscure scientific terminology'
$ws.Range("D5").Value = 'Scientific terms that are simply hard for students to remember.'
$ws.Range("E5").Value = 'Terminology'
$ws.Range("F5").Value = 'Problems with use of language and scientific terms, inconsistent and overlapping terminology.'
$ws.Rows.Item(5).RowHeight = 45

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 'es'
$ws.Range("C6").Value = 'Understanding of Scientific method, process and practice'
$ws.Range("D6").Value = 'Simplistic understandings that may need to be unlearned or revised e.g. imagining atomic structure as balls on sticks suggests space between atoms. Belief that only 50% of parent DNA is passed on to a child. Previous knowledge schemes that need to be modified to integrate new knowledge.'
$ws.Range("E6").Value = 'Incomplete pre-knowledge'
$ws.Range("F6").Value = 'Previous understandings that need to be unlearned, modified or improved to understand the Tricky Topic'
$ws.Rows.Item(6).RowHeight = 120

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 'es'
$ws.Range("C7").Value = 'Underpinning understandings'
$ws.Range("D7").Value = 'Understanding that the student is expected to know already. e.g. to do the calculations related to Avogadro’s number in Chemistry assumes a math understanding of powers of ten and ratios. Learning about genetic drift assumes an understanding of natural selection.'
$ws.Range("E7").Value = 'Incomplete pre-knowledge'
$ws.Range("F7").Value = 'Previous understandings that need to be unlearned, modified or improved to understand the Tricky Topic'
$ws.Rows.Item(7).RowHeight = 120

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'es'
$ws.Range("C8").Value = 'Essential Concepts'
$ws.Range("D8").Value = 'Complementary concepts that the student needs to learn alongside the topic in order to make sense of the new knowledge. e.g. understanding genetic drift involves learning about its causes; founder effect and bottleneck effect.'
$ws.Range("E8").Value = 'Complementary concepts'
$ws.Range("F8").Value = 'Key assumptions and knowledge that relate to the tricky topic, without which it is impossible to understand it'
$ws.Rows.Item(8).RowHeight = 90

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 'es'
$ws.Range("C9").Value = 'Weak human-like or world-like analogy'
$ws.Range("D9").Value = 'Human-Like or world like analogy. Viewing scientific concepts in terms of everyday phenomena e.g. males of any species are bigger than females. Plants suck up food from soil thru roots.Analogy based on metaphor that doesn’t carry through e.g. “Stage” and “Costume” used in Sense programming.'
$ws.Range("E9").Value = 'Intuitive Beliefs'
$ws.Range("F9").Value = 'Informal, intuitive ways of thinking about the world. Strongly biased toward causal explanations'
$ws.Rows.Item(9).RowHeight = 120

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 'es'
$ws.Range("C10").Value = 'Key characteristic conveys group membership'
$ws.Range("D10").Value = 'The belief that if one condition is fulfilled, then the object is automatically a member of a groupOne unobservable core feature defines membership of a category eg: one to one relationship between DNA and physical traits. Birds have wings therefore all creatures with wings are birds.'
$ws.Range("E10").Value = 'Intuitive Beliefs'
$ws.Range("F10").Value = 'Informal, intuitive ways of thinking about the world. Strongly biased toward causal explanations'
$ws.Rows.Item(10).RowHeight = 120

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 'es'
$ws.Range("C11").Value = 'Flawed causal reasoning'
$ws.Range("D11").Value = 'Reasononing based on the assumption of goal or purpose eg birds have wings so they can fly. Genes turn off in order to enable a cell to develop properly. Inappropriate assumption of cause and effect, eg release an object along a curved path and it will continue in a curve, rocks are pointy so that animals won’t sit on them and crush them.'
$ws.Range("E11").Value = 'Intuitive Beliefs'
$ws.Range("F11").Value = 'Informal, intuitive ways of thinking about the world. Strongly biased toward causal explanations'
$ws.Rows.Item(11).RowHeight = 135

# Update the selection / view to match the new table extent
$null = $excel.Goto($ws.Range("B2:B11"), $true)
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
